$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" both contain the same table that needs updating.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Row 2: F2 6003 -> 6011, G2 55 -> 58
    $ws.Range("F2").Value = 6011
    $ws.Range("G2").Value = 58

    # Row 5: F5 993 -> 995
    $ws.Range("F5").Value = 995

    # Row 6: F6 88 -> 89
    $ws.Range("F6").Value = 89
}
